$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5683289
$ws.Range("I33").Value = 1756.8462
$ws.Range("K33").Value = 1756.8462
$ws.Range("M33").Value = -1527.8462

$ws.Range("H101").Value = 1661
$ws.Range("J101").Value = 3491.4285
$ws.Range("L101").Value = 10474.2855
$ws.Range("N101").Value = -13718.2855

$ws.Range("H129").Value = 964.81335
$ws.Range("J129").Value = 985.8472
$ws.Range("L129").Value = 2957.5416
$ws.Range("N129").Value = -12957.5416

$ws.Range("H132").Value = 3054.4119
$ws.Range("I132").Value = 3054.4119
$ws.Range("K132").Value = 9163.235700000001
$ws.Range("M132").Value = -6633.235700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4130.905
$ws.Range("I61").Value = 5450.207
$ws.Range("J61").Value = 1187.8462
$ws.Range("K61").Value = 5450.207
$ws.Range("L61").Value = 1187.8462
$ws.Range("M61").Value = -5238.207
$ws.Range("N61").Value = -1611.8462

$ws.Range("H114").Value = 200000
$ws.Range("J114").Value = 200000
$ws.Range("L114").Value = 200000
$ws.Range("N114").Value = -208678

$ws.Range("H136").Value = 4130.905
$ws.Range("I136").Value = 5450.207
$ws.Range("J136").Value = 1187.8462
$ws.Range("K136").Value = 16350.621
$ws.Range("L136").Value = 3563.5386
$ws.Range("M136").Value = -13800.621
$ws.Range("N136").Value = -8663.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3864.6667
$ws.Range("I134").Value = 5189.5557
$ws.Range("J134").Value = 2374.1667
$ws.Range("K134").Value = 15568.6671
$ws.Range("L134").Value = 7122.500100000001
$ws.Range("M134").Value = -13033.6671
$ws.Range("N134").Value = -12192.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2532.6914
$ws.Range("I31").Value = 2219.611
$ws.Range("J31").Value = 2622.1428
$ws.Range("K31").Value = 2219.611
$ws.Range("L31").Value = 2622.1428
$ws.Range("M31").Value = -1924.611
$ws.Range("N31").Value = -3212.1428

$ws.Range("H34").Value = 2532.6914
$ws.Range("I34").Value = 2219.611
$ws.Range("J34").Value = 2622.1428
$ws.Range("K34").Value = 2219.611
$ws.Range("L34").Value = 2622.1428
$ws.Range("M34").Value = -2017.611
$ws.Range("N34").Value = -3026.1428

$ws.Range("H105").Value = 37039140
$ws.Range("I105").Value = 47620892
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 47620892
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -47619145
$ws.Range("N105").Value = -6494

$ws.Range("H122").Value = 1774.625
$ws.Range("I122").Value = 1867.1666
$ws.Range("J122").Value = 1497
$ws.Range("K122").Value = 5601.4998
$ws.Range("L122").Value = 4491
$ws.Range("M122").Value = -3151.4998
$ws.Range("N122").Value = -9391

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 6268824
$ws.Range("I33").Value = 11111182
$ws.Range("J33").Value = 42935
$ws.Range("K33").Value = 66667092
$ws.Range("L33").Value = 257610
$ws.Range("M33").Value = -66666809
$ws.Range("N33").Value = -258176

$ws.Range("H63").Value = 3750
$ws.Range("I63").Value = 1400
$ws.Range("J63").Value = 6100
$ws.Range("K63").Value = 4200
$ws.Range("L63").Value = 18300
$ws.Range("M63").Value = -3451
$ws.Range("N63").Value = -19798

$ws.Range("H66").Value = 3750
$ws.Range("I66").Value = 1400
$ws.Range("J66").Value = 6100
$ws.Range("K66").Value = 12600
$ws.Range("L66").Value = 54900
$ws.Range("M66").Value = -8856
$ws.Range("N66").Value = -62388

$ws.Range("H68").Value = 2706.935
$ws.Range("I68").Value = 3419.3333
$ws.Range("J68").Value = 1975.7894
$ws.Range("K68").Value = 10257.9999
$ws.Range("L68").Value = 5927.3682
$ws.Range("M68").Value = -9446.999899999999
$ws.Range("N68").Value = -7549.3682

$ws.Range("H71").Value = 2706.935
$ws.Range("I71").Value = 3419.3333
$ws.Range("J71").Value = 1975.7894
$ws.Range("K71").Value = 30773.9997
$ws.Range("L71").Value = 17782.1046
$ws.Range("M71").Value = -26717.9997
$ws.Range("N71").Value = -25894.1046

$ws.Range("H81").Value = 750
$ws.Range("I81").Value = 750
$ws.Range("K81").Value = 2250
$ws.Range("M81").Value = -1127

$ws.Range("H84").Value = 750
$ws.Range("I84").Value = 750
$ws.Range("K84").Value = 6750
$ws.Range("M84").Value = -1134

$ws.Range("H107").Value = 1141.7307
$ws.Range("I107").Value = 333
$ws.Range("J107").Value = 1358.7073
$ws.Range("K107").Value = 999
$ws.Range("L107").Value = 4076.1219
$ws.Range("M107").Value = 921
$ws.Range("N107").Value = -7916.1219

$ws.Range("H108").Value = 100200
$ws.Range("I108").Value = 100200
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 300600
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -297720
$ws.Range("N108").ClearContents()

$ws.Range("H122").Value = 577.8333
$ws.Range("I122").Value = 474.36365
$ws.Range("J122").Value = 623.36
$ws.Range("K122").Value = 4269.27285
$ws.Range("L122").Value = 5610.24
$ws.Range("M122").Value = -1819.27285
$ws.Range("N122").Value = -10510.24

$ws.Range("H131").Value = 11957756
$ws.Range("I131").Value = 5882832.5
$ws.Range("J131").Value = 13334738
$ws.Range("K131").Value = 17648497.5
$ws.Range("L131").Value = 40004214
$ws.Range("M131").Value = -17643457.5
$ws.Range("N131").Value = -40014294

$ws.Range("H132").Value = 2425.8076
$ws.Range("J132").Value = 2501.5352
$ws.Range("L132").Value = 22513.8168
$ws.Range("N132").Value = -27573.8168

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5964.525
$ws.Range("I70").Value = 5939.7334
$ws.Range("J70").Value = 6038.9
$ws.Range("K70").Value = 5939.7334
$ws.Range("L70").Value = 6038.9
$ws.Range("M70").Value = -5669.7334
$ws.Range("N70").Value = -6578.9

$ws.Range("H73").Value = 5964.525
$ws.Range("I73").Value = 5939.7334
$ws.Range("J73").Value = 6038.9
$ws.Range("K73").Value = 5939.7334
$ws.Range("L73").Value = 6038.9
$ws.Range("M73").Value = -5003.7334
$ws.Range("N73").Value = -7910.9

$ws.Range("H126").Value = 5011.0347
$ws.Range("I126").Value = 5678.75
$ws.Range("J126").Value = 1806
$ws.Range("K126").Value = 17036.25
$ws.Range("L126").Value = 5418
$ws.Range("M126").Value = -14566.25
$ws.Range("N126").Value = -10358

$ws.Range("H132").Value = 3634.327
$ws.Range("I132").Value = 4523.923
$ws.Range("J132").Value = 3337.795
$ws.Range("K132").Value = 13571.769
$ws.Range("L132").Value = 10013.385
$ws.Range("M132").Value = -11041.769
$ws.Range("N132").Value = -15073.385

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 166669840
$ws.Range("I40").Value = 200002800
$ws.Range("J40").Value = 4999
$ws.Range("K40").Value = 200002800
$ws.Range("L40").Value = 4999
$ws.Range("M40").Value = -200002664
$ws.Range("N40").Value = -5271

$ws.Range("H61").Value = 2230.7646
$ws.Range("I61").Value = 1708.4286
$ws.Range("J61").Value = 4668.3335
$ws.Range("K61").Value = 1708.4286
$ws.Range("L61").Value = 4668.3335
$ws.Range("M61").Value = -1506.4286
$ws.Range("N61").Value = -5072.3335

$ws.Range("H113").Value = 2230.7646
$ws.Range("I113").Value = 1708.4286
$ws.Range("J113").Value = 4668.3335
$ws.Range("K113").Value = 1708.4286
$ws.Range("L113").Value = 4668.3335
$ws.Range("M113").Value = 461.5714
$ws.Range("N113").Value = -9008.333500000001

$ws.Range("H136").Value = 6547.7114
$ws.Range("I136").Value = 4970.343
$ws.Range("J136").Value = 9795.235000000001
$ws.Range("K136").Value = 14911.029
$ws.Range("L136").Value = 29385.705
$ws.Range("M136").Value = -12361.029
$ws.Range("N136").Value = -34485.705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 26739.4
$ws.Range("J123").Value = 26739.4
$ws.Range("L123").Value = 26739.4
$ws.Range("N123").Value = -36539.4

$ws.Range("H132").Value = 2392.0286
$ws.Range("I132").Value = 1928.2354
$ws.Range("J132").Value = 2830.0557
$ws.Range("K132").Value = 5784.706200000001
$ws.Range("L132").Value = 8490.167099999999
$ws.Range("M132").Value = -3254.706200000001
$ws.Range("N132").Value = -13550.1671
